# Lab 1 (Beverton-Holt) deck update:
#  - Insert a new slide at position 4 containing R code (library/compile/dyn.load/
#    MakeADFun/... ) that walks through setting up and calling the TMB objective,
#    pushing the existing "Exercises #1", "Exercises #2" and "Recap" slides down
#    by one position each.

$p = $ppt.ActivePresentation

# Slide 4 ("Beverton-Hold Exercises #1") currently uses the "Title and Content"
# custom layout; reuse the same layout for the new slide so placeholder
# inheritance (fonts, bullet defaults, etc.) matches.
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

$new = $p.Slides.AddSlide(4, $titleAndContent)

# This slide has no title -- drop the title placeholder entirely.
$new.Shapes.Item(1).Delete()

$content = $new.Shapes.Item(1)

# Reposition/resize the content placeholder now that it's the only shape on
# the slide (it previously sat below a title).
$content.Left = 49.5
$content.Top = 54.4186
$content.Width = 621.0
$content.Height = 431.9565
$content.TextFrame.AutoSize = 2

$lines = @(
    'library(TMB)',
    'compile("tmb_models/bevholt.cpp")',
    'dyn.load(dynlib("tmb_models/bevholt"))',
    '',
    'dat <- read.table("tmb_models/bevholt.dat", header=TRUE)',
    'data <- list(SSB=dat$ssb,logR=dat$logR)',
    'parameters <- list(logA=0, logB=0)',
    '',
    'obj <- MakeADFun(data,parameters,DLL="bevholt")',
    'obj$env$beSilent() # silences console output',
    'obj$fn()',
    'obj$gr()'
)

$content.TextFrame.TextRange.Text = $lines -join "`r"

$tr = $content.TextFrame.TextRange
for ($i = 1; $i -le $lines.Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = $false
    $para.Font.Name = "Courier New"
    $para.Font.Size = 20
}
